$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 741, shifting existing rows 741-758 down to 742-759.
$ws.Rows.Item(741).Insert()

# Populate the newly inserted row 741 with the new weekly price record.
$ws.Cells.Item(741, 1).Value = 6
$ws.Cells.Item(741, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(741, 3).Value = "Metropolitana"
$ws.Cells.Item(741, 4).Value = 44628
$ws.Cells.Item(741, 5).Value = 13
$ws.Cells.Item(741, 6).Value = 100112021
$ws.Cells.Item(741, 7).Value = "Ají"
$ws.Cells.Item(741, 8).Value = "Americana (o)"
$ws.Cells.Item(741, 9).Value = "Primera"
$ws.Cells.Item(741, 10).Value = 140
$ws.Cells.Item(741, 11).Value = 18000
$ws.Cells.Item(741, 12).Value = 20000
$ws.Cells.Item(741, 13).Value = 19143
$ws.Cells.Item(741, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(741, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(741, 16).Value = 766
$ws.Cells.Item(741, 17).Value = 25
$ws.Cells.Item(741, 18).Value = "Hortaliza"
